$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original Row 4 values (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)
$ws.Range("D2").Value = 44980
$ws.Range("M2").Value = 50

# Row 3 <- original Row 5 values
$ws.Range("D3").Value = 44971
$ws.Range("M3").Value = 25

# Row 4 <- original Row 3 values
$ws.Range("D4").Value = 44973
$ws.Range("M4").Value = 55
$ws.Range("N4").Value = 28000
$ws.Range("O4").Value = 28000
$ws.Range("P4").Value = 28000
$ws.Range("S4").Value = 3500

# Row 5 <- original Row 2 values
$ws.Range("D5").Value = 44981
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 25000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 25000
$ws.Range("S5").Value = 3125
